$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

Set-TextValue 'D2' '28.241.79'
$ws.Range('E2').Value = '  -0.67%  '
Set-TextValue 'D3' '1.833.01'
$ws.Range('E3').Value = '  +1.15%  '
$ws.Range('E4').Value = '  +0.03%  '
Set-TextValue 'D5' '310.49'
$ws.Range('E5').Value = '  -0.76%  '
Set-TextValue 'D6' '1.000'
$ws.Range('E6').Value = '  +0.03%  '
Set-TextValue 'D7' '0.4970'
$ws.Range('E7').Value = '  -3.73%  '
$ws.Range('B8').Value = 'Dogecoin'
$ws.Range('C8').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
Set-TextValue 'D8' '0.1010'
$ws.Range('E8').Value = '  +27.78%  '
$ws.Range('B9').Value = 'Cardano'
$ws.Range('C9').Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
Set-TextValue 'D9' '0.3931'
$ws.Range('E9').Value = '  -1.53%  '
Set-TextValue 'D10' '1.113'
$ws.Range('E10').Value = '  -0.30%  '
Set-TextValue 'D11' '41.21'
$ws.Range('E11').Value = '  +0.72%  '
Set-TextValue 'D12' '6.445'
$ws.Range('E12').Value = '  +0.91%  '
$ws.Range('E13').Value = '  +1.22%  '
$ws.Range('E14').Value = '  +0.08%  '
Set-TextValue 'D15' '1.827.36'
$ws.Range('E15').Value = '  +1.01%  '
Set-TextValue 'D16' '7.344'
$ws.Range('E16').Value = '  +0.03%  '
$ws.Range('E17').Value = '  +5.57%  '
Set-TextValue 'D18' '93.17'
$ws.Range('E18').Value = '  +0.37%  '
Set-TextValue 'D19' '0.06657'
$ws.Range('E19').Value = '  +1.33%  '
Set-TextValue 'D20' '0.9997'
$ws.Range('E20').Value = '  +0.02%  '
Set-TextValue 'D21' '17.29'
$ws.Range('E21').Value = '  -0.48%  '
Set-TextValue 'D22' '6.029'
$ws.Range('E22').Value = '  +0.07%  '
Set-TextValue 'D23' '28.285.21'
$ws.Range('E23').Value = '  -0.63%  '
Set-TextValue 'D24' '11.34'
$ws.Range('E24').Value = '  +1.68%  '
Set-TextValue 'D25' '2.230'
$ws.Range('E25').Value = '  +0.03%  '
Set-TextValue 'D26' '158.20'
$ws.Range('E26').Value = '  -1.88%  '
$ws.Range('E27').Value = '  +1.17%  '
Set-TextValue 'D28' '2.039.58'
$ws.Range('E28').Value = '  +0.94%  '
$ws.Range('E29').Value = '  +1.51%  '
Set-TextValue 'D30' '127.07'
$ws.Range('E30').Value = '  -1.04%  '
Set-TextValue 'D31' '0.1055'
$ws.Range('E31').Value = '  -3.68%  '
Set-TextValue 'D32' '1.043'
$ws.Range('E32').Value = '  -2.97%  '
Set-TextValue 'D33' '5.612'
$ws.Range('E33').Value = '  +0.40%  '
Set-TextValue 'D34' '3.601'
$ws.Range('E34').Value = '  -1.75%  '
$ws.Range('E35').Value = '  -6.34%  '
Set-TextValue 'D36' '9.059'
$ws.Range('E36').Value = '  -1.42%  '
Set-TextValue 'D37' '0.02360'
$ws.Range('E37').Value = '  +0.65%  '
Set-TextValue 'D38' '0.2156'
$ws.Range('E38').Value = '  -1.10%  '
Set-TextValue 'D39' '11.50'
$ws.Range('E39').Value = '  -1.42%  '
Set-TextValue 'D40' '4.995'
Set-TextValue 'D41' '0.6240'
$ws.Range('E41').Value = '  +0.36%  '
$ws.Range('E42').Value = '  +1.78%  '
Set-TextValue 'D43' '0.9995'
$ws.Range('E43').Value = '  +0.03%  '
Set-TextValue 'D44' '13.22'
$ws.Range('E44').Value = '  -0.48%  '
Set-TextValue 'D45' '0.5959'
$ws.Range('E45').Value = '  -0.87%  '
Set-TextValue 'D46' '3.691'
$ws.Range('E46').Value = '  -1.17%  '
Set-TextValue 'D47' '1.269'
$ws.Range('E47').Value = '  -3.25%  '
Set-TextValue 'D48' '124.19'
$ws.Range('E48').Value = '  -1.36%  '
Set-TextValue 'D49' '1.956'
$ws.Range('E49').Value = '  +1.08%  '
Set-TextValue 'D50' '1.183'
$ws.Range('E50').Value = '  -3.34%  '
Set-TextValue 'D51' '1.126'
$ws.Range('E51').Value = '  +4.87%  '
